$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '24.491.03'
$ws.Range('E2').Value = '  +10.54%  '

# Row 3
$ws.Range('D3').Value = '1.680.37'
$ws.Range('E3').Value = '  +5.80%  '

# Row 4
$ws.Range('D4').Value = "'" + '1.002'
$ws.Range('E4').Value = '  -0.36%  '

# Row 5
$ws.Range('D5').Value = "'" + '305.85'
$ws.Range('E5').Value = '  +2.74%  '

# Row 6
$ws.Range('D6').Value = "'" + '0.9969'
$ws.Range('E6').Value = '  +0.68%  '

# Row 7
$ws.Range('D7').Value = "'" + '0.3687'
$ws.Range('E7').Value = '  +1.71%  '

# Row 8
$ws.Range('D8').Value = "'" + '0.3424'
$ws.Range('E8').Value = '  +2.46%  '

# Row 9
$ws.Range('D9').Value = "'" + '47.98'
$ws.Range('E9').Value = '  +16.20%  '

# Row 10
$ws.Range('D10').Value = "'" + '1.162'
$ws.Range('E10').Value = '  +4.10%  '

# Row 11
$ws.Range('D11').Value = "'" + '0.07225'
$ws.Range('E11').Value = '  +4.08%  '

# Row 12
$ws.Range('D12').Value = "'" + '0.9985'
$ws.Range('E12').Value = '  -0.42%  '

# Row 13
$ws.Range('D13').Value = "'" + '6.112'
$ws.Range('E13').Value = '  +5.14%  '

# Row 14
$ws.Range('D14').Value = "'" + '20.17'
$ws.Range('E14').Value = '  +4.08%  '

# Row 15
$ws.Range('D15').Value = "'" + '6.715'
$ws.Range('E15').Value = '  +2.94%  '

# Row 16
$ws.Range('D16').Value = '1.677.26'
$ws.Range('E16').Value = '  +5.40%  '

# Row 17
$ws.Range('D17').Value = "'" + '0.00001100'
$ws.Range('E17').Value = '  +3.64%  '

# Row 18
$ws.Range('D18').Value = "'" + '0.9967'
$ws.Range('E18').Value = '  +0.51%  '

# Row 19
$ws.Range('D19').Value = "'" + '0.06650'
$ws.Range('E19').Value = '  +1.44%  '

# Row 20
$ws.Range('D20').Value = "'" + '80.70'
$ws.Range('E20').Value = '  +6.36%  '

# Row 21
$ws.Range('D21').Value = "'" + '16.44'
$ws.Range('E21').Value = '  +4.04%  '

# Row 22
$ws.Range('D22').Value = "'" + '6.082'
$ws.Range('E22').Value = '  +2.91%  '

# Row 23
$ws.Range('D23').Value = "'" + '12.10'
$ws.Range('E23').Value = '  +4.26%  '

# Row 24
$ws.Range('D24').Value = '24.416.24'
$ws.Range('E24').Value = '  +10.11%  '

# Row 25
$ws.Range('D25').Value = "'" + '2.417'
$ws.Range('E25').Value = '  +1.73%  '

# Row 26
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = "'" + '2.653'
$ws.Range('E26').Value = '  +6.55%  '

# Row 27
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = "'" + '153.05'
$ws.Range('E27').Value = '  +3.35%  '

# Row 28
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'" + '19.42'
$ws.Range('E28').Value = '  +1.54%  '

# Row 29
$ws.Range('B29').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C29').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D29').Value = '1.862.49'
$ws.Range('E29').Value = '  +5.77%  '

# Row 30
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = "'" + '127.16'
$ws.Range('E30').Value = '  +4.79%  '

# Row 31
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = "'" + '6.254'
$ws.Range('E31').Value = '  +6.35%  '

# Row 32
$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').Value = "'" + '4.024'
$ws.Range('E32').Value = '  +1.26%  '

# Row 33
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = "'" + '0.9797'
$ws.Range('E33').Value = '  +6.88%  '

# Row 34
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').Value = "'" + '0.08418'
$ws.Range('E34').Value = '  +3.18%  '

# Row 35
$ws.Range('B35').Value = 'WEMIXTOKEN'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = "'" + '1.694'
$ws.Range('E35').Value = '  +4.93%  '

# Row 36
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').Value = "'" + '12.35'
$ws.Range('E36').Value = '  +5.96%  '

# Row 37
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = "'" + '0.06372'
$ws.Range('E37').Value = '  +6.17%  '

# Row 38
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = "'" + '5.317'
$ws.Range('E38').Value = '  +4.22%  '

# Row 39
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = "'" + '8.674'
$ws.Range('E39').Value = '  +4.31%  '

# Row 40
$ws.Range('D40').Value = "'" + '0.02310'
$ws.Range('E40').Value = '  +6.31%  '

# Row 41
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = "'" + '1.246'
$ws.Range('E41').Value = '  +0.94%  '

# Row 42
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = "'" + '0.2089'
$ws.Range('E42').Value = '  +5.60%  '

# Row 43
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = "'" + '0.6092'
$ws.Range('E43').Value = '  +5.70%  '

# Row 44
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = "'" + '0.9970'
$ws.Range('E44').Value = '  +0.54%  '

# Row 45
$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D45').Value = "'" + '3.762'
$ws.Range('E45').Value = '  +0.09%  '

# Row 46
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'" + '12.90'
$ws.Range('E46').Value = '  -0.32%  '

# Row 47
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = "'" + '0.5870'
$ws.Range('E47').Value = '  +5.79%  '

# Row 48
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = "'" + '125.53'
$ws.Range('E48').Value = '  +0.47%  '

# Row 49
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = "'" + '2.009'
$ws.Range('E49').Value = '  +3.79%  '

# Row 50
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'" + '0.07166'
$ws.Range('E50').Value = '  +6.90%  '

# Row 51
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = "'" + '75.66'
$ws.Range('E51').Value = '  +4.55%  '
